$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct some values in the V_matrix table (bug fixing on strategies)
$ws.Range("F2").Value = 999
$ws.Range("G2").Value = 2001
$ws.Range("F5").Value = 1665
$ws.Range("F8").Value = 1265
$ws.Range("G9").Value = 2001
$ws.Range("E10").Value = 2498

# Update selection to match the last-edited cell
$ws.Range("F5").Select()
